$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the cell values that point at the old "test" / legacy environment
#    so they point at the new sandbox environment instead.
$ws.Range("M2").Value = "https://mirandakate.cabisandbox.com"
$ws.Range("A2").Value = "https://sandbox.cabiclio.com/backoffice/control/main"
$ws.Range("G2").Value = "https://sandbox.cabiclio.com/cabicentral/control/main"
$ws.Range("J2").Value = "https://sandbox.cabiclio.com/warehouse/control/main"

# 2. Rebuild the hyperlinks collection. This engine's Hyperlinks.Delete() only
#    works at the worksheet-collection level (per-item Delete() is a no-op), so
#    every hyperlink has to be deleted and re-created together; the three links
#    whose target text changed (G2, J2, A2) get no explicit display text so
#    Excel just shows the (new) cell text with no stale "display" override,
#    while the untouched ones are recreated with their original target/text.
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("F2"), 'mailto:C@bi$ush5', [Type]::Missing, [Type]::Missing, 'C@bi$ush5')
$ws.Hyperlinks.Add($ws.Range("G2"), "https://test17.cliotest.com/cabicentral/control/main")
$ws.Hyperlinks.Add($ws.Range("J2"), "https://test19.cliotest.com/warehouse/control/main")
$ws.Hyperlinks.Add($ws.Range("A2"), "https://test4.cliotest.com/backoffice/control/main")
$ws.Hyperlinks.Add($ws.Range("M2"), "https://mirandakate.cabitest5.com/")
$ws.Hyperlinks.Add($ws.Range("N2"), "mailto:michigan@na.com")
$ws.Hyperlinks.Add($ws.Range("P2"), "https://mirandakate.cabionline.com/")

# 3. Hyperlinks.Add() re-stamps the "Hyperlink" cell style (and nudges the
#    cell text for ranges that didn't get an explicit display string above),
#    so restore the original cell text/formatting for every touched cell.
$ws.Range("F2").Value = "cabiautomation"
$ws.Range("F2").Font.Name = "Arial"
$ws.Range("G2").Font.Name = "Arial"
$ws.Range("J2").Font.Name = "Arial"
$ws.Range("A2").Font.Name = "Arial"
$ws.Range("M2").Font.Name = "Arial"
$ws.Range("N2").Font.Name = "Arial"
$ws.Range("P2").Font.Name = "Arial"

# 4. Update the active selection, as recorded in the saved view state.
$ws.Range("R2").Select()
